# TC47_Canine_Filter_Breed-YorkshireTerr.xlsx
# "updated cart validation -1 web and db"
#
# Inserts a new "cartQuery" column (D) between the existing dbExcel
# (StatQuery) column and the FilesTab column, containing the new
# cart-validation Cypher query used by the web/db cart export, and
# shifts the old FilesTab / WebExcel columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert the new column D (pushes old D/E -> E/F) -----------------
$ws.Columns.Item(4).Insert()

# --- header row --------------------------------------------------------
$ws.Range("D1").Value = "cartQuery"

# --- new cart query text (same value repeated down col D for rows 2-4) -
$nl = [char]10
$cartQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)" + $nl + `
  "  WHERE demo.breed IN ['Yorkshire Terrier']" + $nl + `
  "MATCH (f:file)-[*]->(c)" + $nl + `
  "WITH COLLECT(DISTINCT f.uuid) AS uuids" + $nl + `
  "MATCH (f:file)" + $nl + `
  "  WHERE f.uuid in uuids" + $nl + `
  "OPTIONAL MATCH (f)-->(parent)" + $nl + `
  "OPTIONAL MATCH (f)-[*]->(samp:sample)" + $nl + `
  "OPTIONAL MATCH (f:file)-[*]->(c:case)" + $nl + `
  "OPTIONAL MATCH (s:study)<-[:member_of]-(c)" + $nl + `
  "OPTIONAL MATCH (c)-->(i:canine_individual)<--(o:case)" + $nl + `
  "RETURN" + $nl + `
  "  f.file_name AS ``File Name``," + $nl + `
  "  f.file_type AS ``File Type``," + $nl + `
  "  head(labels(parent)) AS ``Association``," + $nl + `
  "  f.file_description AS ``Description``," + $nl + `
  "  f.file_format AS ``Format``," + $nl + `
  "  f.file_size AS ``Size``," + $nl + `
  "  samp.sample_id AS ``Sample ID``," + $nl + `
  "  c.case_id as ``Case ID``," + $nl + `
  "  i.canine_individual_id AS ``Canine ID``," + $nl + `
  "  CASE WHEN s.clinical_study_designation IS NULL " + $nl + `
  "  THEN parent.clinical_study_designation " + $nl + `
  "  ELSE s.clinical_study_designation END AS ``Study Code``" + $nl + `
  "  "

$ws.Range("D2").Value = $cartQuery
$ws.Range("D3").Value = $cartQuery
$ws.Range("D4").Value = $cartQuery

# match the wrap-text style used by the other long-text columns (B, C)
$ws.Range("D2:D4").WrapText = $true

# --- row heights grow to fit the longer wrapped query text -------------
$ws.Rows.Item(2).RowHeight = 390
$ws.Rows.Item(3).RowHeight = 390
$ws.Rows.Item(4).RowHeight = 390

# --- restore cursor/selection to where the editor left off -------------
$ws.Range("C14").Select()
